# Update cryptocurrency price/volume figures per the Sun Jun 30 11:09:14 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.556.77"
$ws.Range("E2").Value = "'  +1.02%  "
$ws.Range("D3").Value = "'3.388.99"
$ws.Range("E3").Value = "'  -0.17%  "
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'575.24"
$ws.Range("E5").Value = "'  +0.44%  "
$ws.Range("D6").Value = "'140.96"
$ws.Range("E6").Value = "'  -1.09%  "
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("D8").Value = "'0.473"
$ws.Range("E8").Value = "'  -0.68%  "
$ws.Range("E9").Value = "'  +1.74%  "
$ws.Range("E10").Value = "'  -1.23%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "'  -2.52%  "
$ws.Range("D12").Value = "'3.964.33"
$ws.Range("E12").Value = "'  -0.22%  "
$ws.Range("E13").Value = "'  +0.22%  "
$ws.Range("D14").Value = "'28.35"
$ws.Range("D15").Value = "'3.393.95"
$ws.Range("E15").Value = "'  -0.16%  "
$ws.Range("E16").Value = "'  -0.51%  "
$ws.Range("D17").Value = "'61.542.32"
$ws.Range("E17").Value = "'  +0.88%  "
$ws.Range("E18").Value = "'  -0.53%  "
$ws.Range("D19").Value = "'13.62"
$ws.Range("E19").Value = "'  -1.78%  "
$ws.Range("D20").Value = "'9.01"
$ws.Range("E20").Value = "'  +0.63%  "
$ws.Range("D21").Value = "'391.55"
$ws.Range("E21").Value = "'  +2.19%  "
$ws.Range("D22").Value = "'74.89"
$ws.Range("E22").Value = "'  +0.76%  "
$ws.Range("D23").Value = "'0.552"
$ws.Range("E23").Value = "'  -1.08%  "
$ws.Range("E24").Value = "'  +0.41%  "
$ws.Range("E25").Value = "'  -4.26%  "
$ws.Range("E26").Value = "'  +7.37%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "'  -0.09%  "
$ws.Range("D28").Value = "'7.26"
$ws.Range("E28").Value = "'  -1.64%  "
$ws.Range("D29").Value = "'8.04"
$ws.Range("E29").Value = "'  +0.32%  "
$ws.Range("E30").Value = "'  -0.90%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "'  -0.64%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "'  +0.00%  "
$ws.Range("D33").Value = "'23.27"
$ws.Range("E33").Value = "'  -0.95%  "
$ws.Range("D34").Value = "'169.54"
$ws.Range("E34").Value = "'  +1.13%  "
$ws.Range("D35").Value = "'6.90"
$ws.Range("E35").Value = "'  -1.41%  "
$ws.Range("D36").Value = "'5.04"
$ws.Range("E36").Value = "'  +0.69%  "
$ws.Range("D37").Value = "'3.421.73"
$ws.Range("E37").Value = "'  -0.10%  "
$ws.Range("D38").Value = "'1.47"
$ws.Range("E38").Value = "'  -1.19%  "
$ws.Range("D39").Value = "'0.0767"
$ws.Range("E39").Value = "'  -0.58%  "
$ws.Range("D40").Value = "'25.98"
$ws.Range("E40").Value = "'  -5.94%  "
$ws.Range("D41").Value = "'0.780"
$ws.Range("E41").Value = "'  -0.17%  "
$ws.Range("D42").Value = "'4.42"
$ws.Range("E42").Value = "'  -0.28%  "
$ws.Range("E43").Value = "'  -1.02%  "
$ws.Range("D44").Value = "'1.16"
$ws.Range("E44").Value = "'  +1.98%  "
$ws.Range("D45").Value = "'2.465.86"
$ws.Range("E45").Value = "'  -0.60%  "
$ws.Range("D46").Value = "'22.83"
$ws.Range("E46").Value = "'  -0.72%  "
$ws.Range("D47").Value = "'6.66"
$ws.Range("E47").Value = "'  -2.21%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "'  -0.10%  "
$ws.Range("D49").Value = "'0.0264"
$ws.Range("E49").Value = "'  -0.90%  "
$ws.Range("E50").Value = "'  -5.04%  "
$ws.Range("E51").Value = "'  -1.67%  "
